$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New comment text in J17 (adds a new shared string + enables wrap so the
# long text displays on multiple lines, matching the style used elsewhere
# in column J / I for similar remark cells).
$ws.Range("J17").Value = "Vu avec Michel et la DRHFPNC, notre calcul est ok et ne doit pas tenir compte ni de l'ACC ni de la BM"
$ws.Range("J17").WrapText = $true

# Mark the rows that still need clarification with a red flag in the new
# column K (no text, just a solid red fill).
$ws.Range("K7").Interior.Color = 255
$ws.Range("K13").Interior.Color = 255
$ws.Range("K18").Interior.Color = 255
$ws.Range("K21").Interior.Color = 255

# Leave the selection on the last flag that was added.
$ws.Range("K21").Select()
